$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$newRow = 88
$rng = $ws.Range("A" + $newRow + ":E" + $newRow)

# Every existing row in this sheet stores Date/Phase/Result/etc. as literal
# text (even the date- and number-looking ones), with no explicit cell
# style. Typing straight into .Value would let Excel auto-convert
# "2025-12-13" / "251213" into a real date/number and would stamp a
# non-default NumberFormat style onto the cells. So: force Text format,
# assign the literal values, then clear the formatting back off again -
# the stored values stay text (Excel doesn't re-parse on a format change)
# while the cells end up with no explicit style, matching the rest of the
# sheet.
$rng.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-12-13"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
$ws.Cells.Item($newRow, 3).Value = "251213"
$ws.Cells.Item($newRow, 4).Value = "8-0-8"
$ws.Cells.Item($newRow, 5).Value = "2025-12-13T21:37:17.613+04:00"

$rng.ClearFormats()
